$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_1_10_A")

# --- Update title / column header text (October -> November) ---
# NumberFormat is toggled to "@" (text) around the assignment so Excel's
# COM layer stores the value as a literal string instead of auto-parsing
# "November 2016" into a date serial, then restored to the original
# "#,##0" numeric format so the cell style/formatting is unchanged.
$ws.Range("A2").Value = "by State, by Sector, November 2016 and 2015 (Thousand Megawatthours)"

$headerCells = @("B6","E6","G6","I6","K6")
foreach ($addr in $headerCells) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = "November 2016"
    $ws.Range($addr).NumberFormat = "#,##0"
}
$headerCells2 = @("C6","F6","H6","J6","L6")
foreach ($addr in $headerCells2) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = "November 2015"
    $ws.Range($addr).NumberFormat = "#,##0"
}

# --- Update data values (rows 7-68) to the new EPM_2016_11 run figures ---
$ws.Range("B7").Value = 359
$ws.Range("C7").Value = 474
$ws.Range("D7").Value = -0.242
$ws.Range("E7").Value = 46
$ws.Range("F7").Value = 70
$ws.Range("H7").Value = 375
$ws.Range("J7").Value = 0.45
$ws.Range("K7").Value = 18
$ws.Range("L7").Value = 29
$ws.Range("H8").Value = 19
$ws.Range("B9").Value = 186
$ws.Range("C9").Value = 224
$ws.Range("D9").Value = -0.167
$ws.Range("G9").Value = 168
$ws.Range("H9").Value = 195
$ws.Range("K9").Value = 18
$ws.Range("L9").Value = 28
$ws.Range("B10").Value = 44
$ws.Range("C10").Value = 66
$ws.Range("D10").Value = -0.342
$ws.Range("F10").Value = 16
$ws.Range("G10").Value = 34
$ws.Range("H10").Value = 49
$ws.Range("J10").Value = 0.45
$ws.Range("L10").Value = 0.38
$ws.Range("B11").Value = 61
$ws.Range("C11").Value = 81
$ws.Range("D11").Value = -0.253
$ws.Range("E11").Value = 17
$ws.Range("F11").Value = 23
$ws.Range("G11").Value = 44
$ws.Range("B13").Value = 57
$ws.Range("C13").Value = 82
$ws.Range("D13").Value = -0.308
$ws.Range("F13").Value = 29
$ws.Range("H13").Value = 53
$ws.Range("B14").Value = 2237
$ws.Range("C14").Value = 2468
$ws.Range("D14").Value = -0.094
$ws.Range("E14").Value = 1937
$ws.Range("F14").Value = 1949
$ws.Range("G14").Value = 296
$ws.Range("H14").Value = 514
$ws.Range("J14").Value = 0.4
$ws.Range("B16").Value = 2119
$ws.Range("C16").Value = 2287
$ws.Range("D16").Value = -0.073
$ws.Range("E16").Value = 1935
$ws.Range("F16").Value = 1944
$ws.Range("G16").Value = 181
$ws.Range("H16").Value = 336
$ws.Range("J16").Value = 0.4
$ws.Range("B17").Value = 117
$ws.Range("C17").Value = 181
$ws.Range("D17").Value = -0.355
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 114
$ws.Range("H17").Value = 176
$ws.Range("B18").Value = 420
$ws.Range("C18").Value = 610
$ws.Range("D18").Value = -0.312
$ws.Range("E18").Value = 370
$ws.Range("F18").Value = 519
$ws.Range("H18").Value = 65
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0.13
$ws.Range("L18").Value = 26
$ws.Range("C19").Value = 10
$ws.Range("F19").Value = 4
$ws.Range("H19").Value = 7
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0.13
$ws.Range("B20").Value = 40
$ws.Range("C20").Value = 52
$ws.Range("D20").Value = -0.242
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 52
$ws.Range("B21").Value = 124
$ws.Range("C21").Value = 135
$ws.Range("D21").Value = -0.077
$ws.Range("E21").Value = 115
$ws.Range("F21").Value = 125
$ws.Range("H21").Value = 7
$ws.Range("L21").Value = 3
$ws.Range("B22").Value = 47
$ws.Range("C22").Value = 78
$ws.Range("D22").Value = -0.4
$ws.Range("E22").Value = 35
$ws.Range("F22").Value = 42
$ws.Range("H22").Value = 36
$ws.Range("B23").Value = 200
$ws.Range("C23").Value = 335
$ws.Range("D23").Value = -0.403
$ws.Range("E23").Value = 177
$ws.Range("F23").Value = 297
$ws.Range("H23").Value = 15
$ws.Range("L23").Value = 23
$ws.Range("B24").Value = 762
$ws.Range("C24").Value = 788
$ws.Range("D24").Value = -0.034
$ws.Range("E24").Value = 734
$ws.Range("F24").Value = 754
$ws.Range("H24").Value = 21
$ws.Range("L24").Value = 13
$ws.Range("B25").Value = 82
$ws.Range("C25").Value = 60
$ws.Range("D25").Value = 0.373
$ws.Range("E25").Value = 82
$ws.Range("F25").Value = 60
$ws.Range("H25").Value = 0.4
$ws.Range("B27").Value = 74
$ws.Range("C27").Value = 92
$ws.Range("D27").Value = -0.2
$ws.Range("E27").Value = 48
$ws.Range("F27").Value = 60
$ws.Range("H27").Value = 20
$ws.Range("L27").Value = 13
$ws.Range("B28").Value = 31
$ws.Range("C28").Value = 63
$ws.Range("D28").Value = -0.51
$ws.Range("E28").Value = 31
$ws.Range("F28").Value = 63
$ws.Range("B29").Value = 144
$ws.Range("C29").Value = 105
$ws.Range("D29").Value = 0.372
$ws.Range("E29").Value = 144
$ws.Range("F29").Value = 105
$ws.Range("B30").Value = 119
$ws.Range("C30").Value = 117
$ws.Range("D30").Value = 0.014
$ws.Range("E30").Value = 119
$ws.Range("F30").Value = 117
$ws.Range("B31").Value = 310
$ws.Range("C31").Value = 349
$ws.Range("D31").Value = -0.112
$ws.Range("E31").Value = 310
$ws.Range("F31").Value = 349
$ws.Range("B32").Value = 635
$ws.Range("C32").Value = 1891
$ws.Range("D32").Value = -0.664
$ws.Range("E32").Value = 514
$ws.Range("F32").Value = 1669
$ws.Range("H32").Value = 170
$ws.Range("J32").Value = 2
$ws.Range("K32").Value = 16
$ws.Range("L32").Value = 51
$ws.Range("C35").Value = 31
$ws.Range("F35").Value = 31
$ws.Range("B36").Value = 166
$ws.Range("C36").Value = 383
$ws.Range("D36").Value = -0.566
$ws.Range("E36").Value = 164
$ws.Range("F36").Value = 379
$ws.Range("L36").Value = 3
$ws.Range("B37").Value = 59
$ws.Range("C37").Value = 106
$ws.Range("D37").Value = -0.445
$ws.Range("G37").Value = 59
$ws.Range("H37").Value = 106
$ws.Range("B38").Value = 160
$ws.Range("C38").Value = 667
$ws.Range("D38").Value = -0.76
$ws.Range("E38").Value = 157
$ws.Range("F38").Value = 659
$ws.Range("H38").Value = 6
$ws.Range("B39").Value = 100
$ws.Range("C39").Value = 517
$ws.Range("D39").Value = -0.807
$ws.Range("E39").Value = 97
$ws.Range("F39").Value = 504
$ws.Range("H39").Value = 13
$ws.Range("J39").Value = 0.47
$ws.Range("B40").Value = 50
$ws.Range("C40").Value = 53
$ws.Range("D40").Value = -0.061
$ws.Range("E40").Value = 45
$ws.Range("F40").Value = 49
$ws.Range("H40").Value = 3
$ws.Range("B41").Value = 87
$ws.Range("C41").Value = 135
$ws.Range("D41").Value = -0.358
$ws.Range("F41").Value = 47
$ws.Range("G41").Value = 37
$ws.Range("H41").Value = 40
$ws.Range("K41").Value = 13
$ws.Range("L41").Value = 48
$ws.Range("B42").Value = 815
$ws.Range("C42").Value = 2221
$ws.Range("D42").Value = -0.633
$ws.Range("E42").Value = 814
$ws.Range("F42").Value = 2220
$ws.Range("B43").Value = 268
$ws.Range("C43").Value = 1024
$ws.Range("D43").Value = -0.738
$ws.Range("E43").Value = 268
$ws.Range("F43").Value = 1024
$ws.Range("B44").Value = 186
$ws.Range("C44").Value = 258
$ws.Range("D44").Value = -0.278
$ws.Range("E44").Value = 186
$ws.Range("F44").Value = 257
$ws.Range("B46").Value = 360
$ws.Range("C46").Value = 939
$ws.Range("D46").Value = -0.616
$ws.Range("E46").Value = 360
$ws.Range("F46").Value = 939
$ws.Range("B47").Value = 271
$ws.Range("C47").Value = 430
$ws.Range("D47").Value = -0.369
$ws.Range("E47").Value = 210
$ws.Range("F47").Value = 380
$ws.Range("G47").Value = 61
$ws.Range("H47").Value = 51
$ws.Range("B48").Value = 114
$ws.Range("C48").Value = 184
$ws.Range("D48").Value = -0.379
$ws.Range("E48").Value = 111
$ws.Range("F48").Value = 181
$ws.Range("B49").Value = 56
$ws.Range("C49").Value = 45
$ws.Range("D49").Value = 0.253
$ws.Range("G49").Value = 56
$ws.Range("H49").Value = 45
$ws.Range("B50").Value = 63
$ws.Range("C50").Value = 136
$ws.Range("D50").Value = -0.538
$ws.Range("E50").Value = 63
$ws.Range("F50").Value = 136
$ws.Range("B51").Value = 38
$ws.Range("C51").Value = 65
$ws.Range("D51").Value = -0.417
$ws.Range("E51").Value = 36
$ws.Range("F51").Value = 62
$ws.Range("H51").Value = 3
$ws.Range("B52").Value = 2522
$ws.Range("C52").Value = 1717
$ws.Range("D52").Value = 0.469
$ws.Range("E52").Value = 2452
$ws.Range("F52").Value = 1682
$ws.Range("G52").Value = 69
$ws.Range("H52").Value = 35
$ws.Range("J52").Value = 0.05
$ws.Range("B53").Value = 467
$ws.Range("C53").Value = 419
$ws.Range("D53").Value = 0.115
$ws.Range("E53").Value = 467
$ws.Range("F53").Value = 419
$ws.Range("B54").Value = 100
$ws.Range("C54").Value = 14
$ws.Range("D54").Value = 6.349
$ws.Range("E54").Value = 83
$ws.Range("F54").Value = 12
$ws.Range("H54").Value = 2
$ws.Range("J54").Value = 0.05
$ws.Range("B55").Value = 599
$ws.Range("C55").Value = 407
$ws.Range("D55").Value = 0.474
$ws.Range("E55").Value = 560
$ws.Range("F55").Value = 384
$ws.Range("G55").Value = "NM"
$ws.Range("H55").Value = 22
$ws.Range("B56").Value = 1064
$ws.Range("C56").Value = 627
$ws.Range("D56").Value = 0.697
$ws.Range("E56").Value = 1052
$ws.Range("F56").Value = 618
$ws.Range("H56").Value = 9
$ws.Range("B57").Value = 184
$ws.Range("C57").Value = 162
$ws.Range("D57").Value = 0.133
$ws.Range("E57").Value = 183
$ws.Range("F57").Value = 161
$ws.Range("B59").Value = 64
$ws.Range("C59").Value = 45
$ws.Range("D59").Value = 0.408
$ws.Range("E59").Value = 63
$ws.Range("F59").Value = 45
$ws.Range("H59").Value = 1
$ws.Range("B60").Value = 34
$ws.Range("C60").Value = 37
$ws.Range("D60").Value = -0.088
$ws.Range("E60").Value = 33
$ws.Range("F60").Value = 37
$ws.Range("H60").Value = 0.43
$ws.Range("B61").Value = 10649
$ws.Range("C61").Value = 8549
$ws.Range("D61").Value = 0.246
$ws.Range("E61").Value = 10526
$ws.Range("F61").Value = 8485
$ws.Range("G61").Value = 121
$ws.Range("H61").Value = 63
$ws.Range("J61").Value = 0.16
$ws.Range("B62").Value = 1674
$ws.Range("C62").Value = 640
$ws.Range("D62").Value = 1.617
$ws.Range("E62").Value = 1601
$ws.Range("F62").Value = 616
$ws.Range("G62").Value = 71
$ws.Range("H62").Value = 24
$ws.Range("J62").Value = 0.16
$ws.Range("B63").Value = 2772
$ws.Range("C63").Value = 2533
$ws.Range("D63").Value = 0.094
$ws.Range("E63").Value = 2749
$ws.Range("F63").Value = 2513
$ws.Range("H63").Value = 20
$ws.Range("B64").Value = 6203
$ws.Range("C64").Value = 5376
$ws.Range("D64").Value = 0.154
$ws.Range("E64").Value = 6176
$ws.Range("F64").Value = 5356
$ws.Range("G64").Value = "NM"
$ws.Range("H64").Value = 20
$ws.Range("B65").Value = 146
$ws.Range("C65").Value = 189
$ws.Range("D65").Value = -0.226
$ws.Range("E65").Value = 137
$ws.Range("F65").Value = 174
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 6
$ws.Range("L65").Value = 9
$ws.Range("B66").Value = 135
$ws.Range("C66").Value = 170
$ws.Range("D66").Value = -0.206
$ws.Range("E66").Value = 135
$ws.Range("F66").Value = 170
$ws.Range("B67").Value = "NM"
$ws.Range("C67").Value = 19
$ws.Range("D67").Value = "NM"
$ws.Range("F67").Value = 3
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 6
$ws.Range("L67").Value = 9
$ws.Range("B68").Value = 18815
$ws.Range("C68").Value = 19338
$ws.Range("D68").Value = -0.027
$ws.Range("E68").Value = 17741
$ws.Range("F68").Value = 17901
$ws.Range("G68").Value = 1003
$ws.Range("H68").Value = 1301
$ws.Range("K68").Value = 68
$ws.Range("L68").Value = 133

Write-Host "Applied EIA Table 1.10.A November 2016/2015 data update."
